{"js": "// Replace the date line and the 25 \"two-digit \u00f7 one-digit\" answer cells\n// with their updated values, while preserving each run's existing\n// formatting (font, size, paragraph alignment, etc.). Every \"old\" string\n// below is unique within the document, so a body-wide search reliably\n// finds the single matching run of text.\nconst replacements = [\n  [\"2025-05-06 Tuesday\", \"2025-05-07 Wednesday\"],\n  [\"24\u00f75=4, 4\", \"77\u00f79=8, 5\"],\n  [\"76\u00f75=15, 1\", \"32\u00f76=5, 2\"],\n  [\"56\u00f76=9, 2\", \"74\u00f75=14, 4\"],\n  [\"64\u00f72=32, 0\", \"12\u00f74=3, 0\"],\n  [\"75\u00f75=15, 0\", \"64\u00f78=8, 0\"],\n  [\"75\u00f73=25, 0\", \"87\u00f74=21, 3\"],\n  [\"61\u00f76=10, 1\", \"75\u00f75=15, 0\"],\n  [\"99\u00f74=24, 3\", \"16\u00f75=3, 1\"],\n  [\"42\u00f76=7, 0\", \"86\u00f73=28, 2\"],\n  [\"21\u00f74=5, 1\", \"66\u00f72=33, 0\"],\n  [\"27\u00f72=13, 1\", \"71\u00f78=8, 7\"],\n  [\"22\u00f74=5, 2\", \"29\u00f79=3, 2\"],\n  [\"11\u00f72=5, 1\", \"96\u00f73=32, 0\"],\n  [\"70\u00f78=8, 6\", \"89\u00f75=17, 4\"],\n  [\"48\u00f78=6, 0\", \"92\u00f77=13, 1\"],\n  [\"44\u00f74=11, 0\", \"74\u00f75=14, 4\"],\n  [\"55\u00f79=6, 1\", \"48\u00f76=8, 0\"],\n  [\"77\u00f75=15, 2\", \"78\u00f76=13, 0\"],\n  [\"20\u00f78=2, 4\", \"21\u00f74=5, 1\"],\n  [\"36\u00f76=6, 0\", \"28\u00f78=3, 4\"],\n  [\"82\u00f75=16, 2\", \"91\u00f73=30, 1\"],\n  [\"41\u00f76=6, 5\", \"62\u00f77=8, 6\"],\n  [\"20\u00f77=2, 6\", \"69\u00f76=11, 3\"],\n  [\"50\u00f74=12, 2\", \"26\u00f72=13, 0\"],\n  [\"23\u00f73=7, 2\", \"74\u00f76=12, 2\"],\n];\n\nconst body = context.document.body;\n\n// First pass: locate every target range (each old string is unique in\n// the document), so later insertions never accidentally match an\n// earlier/later replacement's \"old\" or \"new\" text.\nconst searchResults = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true })\n);\nsearchResults.forEach((results) => results.load(\"items\"));\nawait context.sync();\n\n// Second pass: replace each found range's text in place, which keeps the\n// surrounding run/paragraph formatting untouched.\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const items = searchResults[i].items;\n  if (items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${replacements[i][0]}\", found ${items.length}`\n    );\n  }\n  items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 \"two-digit \u00f7 one-digit\" answer cells\n# with new values. Assigning directly to Range.Text replaces only the\n# visible text of the paragraph/cell and leaves the existing run/paragraph\n# formatting (fonts, size, alignment, etc.) untouched.\n$d = $word.ActiveDocument\n\n# Heading paragraph with the date.\n$d.Paragraphs.Item(1).Range.Text = \"2025-05-07 Wednesday\"\n\n$t = $d.Tables.Item(1)\n\n# New values for each data row (table rows 1, 5, 9, 13, 17 hold the\n# answers; the rows in between are blank spacer rows), 5 columns each.\n$newValues = @(\n  @(\"77\u00f79=8, 5\", \"32\u00f76=5, 2\", \"74\u00f75=14, 4\", \"12\u00f74=3, 0\", \"64\u00f78=8, 0\"),\n  @(\"87\u00f74=21, 3\", \"75\u00f75=15, 0\", \"16\u00f75=3, 1\", \"86\u00f73=28, 2\", \"66\u00f72=33, 0\"),\n  @(\"71\u00f78=8, 7\", \"29\u00f79=3, 2\", \"96\u00f73=32, 0\", \"89\u00f75=17, 4\", \"92\u00f77=13, 1\"),\n  @(\"74\u00f75=14, 4\", \"48\u00f76=8, 0\", \"78\u00f76=13, 0\", \"21\u00f74=5, 1\", \"28\u00f78=3, 4\"),\n  @(\"91\u00f73=30, 1\", \"62\u00f77=8, 6\", \"69\u00f76=11, 3\", \"26\u00f72=13, 0\", \"74\u00f76=12, 2\")\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n  $row = $dataRows[$i]\n  for ($col = 1; $col -le 5; $col++) {\n    $t.Cell($row, $col).Range.Text = $newValues[$i][$col - 1]\n  }\n}\n"}
